# Registers.xlsx update: add TMF8828 register sheet (copied/adapted from the
# TMF8801 sheet) and remove the now-redundant "8828"-prefixed rows that used
# to live at the bottom of the TMF8801 sheet.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "TMF8828" worksheet right after "TMF8801".
# ------------------------------------------------------------------
$tmf8801 = $wb.Worksheets.Item("TMF8801")
$newSheet = $wb.Worksheets.Add($null, $tmf8801)
$newSheet.Name = "TMF8828"

# ------------------------------------------------------------------
# 2. Fill in the register table for TMF8828.
# ------------------------------------------------------------------
$newSheet.Cells.Item(1,1).Value = 'Name'
$newSheet.Cells.Item(1,2).Value = 'Hex Address'
$newSheet.Cells.Item(1,3).Value = 'Default Value'
$newSheet.Cells.Item(1,4).Value = 'Bit Width'
$newSheet.Cells.Item(1,5).Value = 'Bit Index (High)'
$newSheet.Cells.Item(1,6).Value = 'Bit Index (Low)'
$newSheet.Cells.Item(2,1).Value = 'APPID'
$newSheet.Cells.Item(2,2).Value = '0X00'
$newSheet.Cells.Item(2,3).Value = '0X00'
$newSheet.Cells.Item(2,4).Value = 8
$newSheet.Cells.Item(2,5).Value = 7
$newSheet.Cells.Item(2,6).Value = 0
$newSheet.Cells.Item(3,1).Value = 'MINOR'
$newSheet.Cells.Item(3,2).Value = '0X01'
$newSheet.Cells.Item(3,3).Value = '0X00'
$newSheet.Cells.Item(3,4).Value = 8
$newSheet.Cells.Item(3,5).Value = 7
$newSheet.Cells.Item(3,6).Value = 0
$newSheet.Cells.Item(4,1).Value = 'ENABLE'
$newSheet.Cells.Item(4,2).Value = '0XE0'
$newSheet.Cells.Item(4,3).Value = '0X01'
$newSheet.Cells.Item(4,4).Value = 8
$newSheet.Cells.Item(4,5).Value = 7
$newSheet.Cells.Item(4,6).Value = 0
$newSheet.Cells.Item(5,1).Value = 'CPU_READY'
$newSheet.Cells.Item(5,2).Value = '0xE0'
$newSheet.Cells.Item(5,3).Value = '0X00'
$newSheet.Cells.Item(5,4).Value = 1
$newSheet.Cells.Item(5,5).Value = 6
$newSheet.Cells.Item(5,6).Value = 6
$newSheet.Cells.Item(6,1).Value = 'PON'
$newSheet.Cells.Item(6,2).Value = '0xE0'
$newSheet.Cells.Item(6,3).Value = '0X00'
$newSheet.Cells.Item(6,4).Value = 1
$newSheet.Cells.Item(6,5).Value = 0
$newSheet.Cells.Item(6,6).Value = 0
$newSheet.Cells.Item(7,1).Value = 'INT_STATUS'
$newSheet.Cells.Item(7,2).Value = '0XE1'
$newSheet.Cells.Item(7,3).Value = '0X00'
$newSheet.Cells.Item(7,4).Value = 8
$newSheet.Cells.Item(7,5).Value = 7
$newSheet.Cells.Item(7,6).Value = 0
$newSheet.Cells.Item(8,1).Value = 'INT_ENAB'
$newSheet.Cells.Item(8,2).Value = '0XE2'
$newSheet.Cells.Item(8,3).Value = '0X00'
$newSheet.Cells.Item(8,4).Value = 8
$newSheet.Cells.Item(8,5).Value = 7
$newSheet.Cells.Item(8,6).Value = 0
$newSheet.Cells.Item(9,1).Value = 'ID'
$newSheet.Cells.Item(9,2).Value = '0XE3'
$newSheet.Cells.Item(9,3).Value = '0X08'
$newSheet.Cells.Item(9,4).Value = 6
$newSheet.Cells.Item(9,5).Value = 5
$newSheet.Cells.Item(9,6).Value = 0
$newSheet.Cells.Item(10,1).Value = 'REVID'
$newSheet.Cells.Item(10,2).Value = '0XE4'
$newSheet.Cells.Item(10,3).Value = '0X00'
$newSheet.Cells.Item(10,4).Value = 3
$newSheet.Cells.Item(10,5).Value = 2
$newSheet.Cells.Item(10,6).Value = 0
$newSheet.Cells.Item(11,1).Value = 'PATCH'
$newSheet.Cells.Item(11,2).Value = '0X02'
$newSheet.Cells.Item(11,3).Value = '0X00'
$newSheet.Cells.Item(11,4).Value = 8
$newSheet.Cells.Item(11,5).Value = 7
$newSheet.Cells.Item(11,6).Value = 0
$newSheet.Cells.Item(12,1).Value = 'BUILD_TYPE'
$newSheet.Cells.Item(12,2).Value = '0X03'
$newSheet.Cells.Item(12,3).Value = '0X00'
$newSheet.Cells.Item(12,4).Value = 8
$newSheet.Cells.Item(12,5).Value = 7
$newSheet.Cells.Item(12,6).Value = 0
$newSheet.Cells.Item(13,1).Value = 'APPLICATION_STATUS'
$newSheet.Cells.Item(13,2).Value = '0X04'
$newSheet.Cells.Item(13,3).Value = '0X00'
$newSheet.Cells.Item(13,4).Value = 8
$newSheet.Cells.Item(13,5).Value = 7
$newSheet.Cells.Item(13,6).Value = 0
$newSheet.Cells.Item(14,1).Value = 'MEASURE_STATUS'
$newSheet.Cells.Item(14,2).Value = '0X05'
$newSheet.Cells.Item(14,3).Value = '0X00'
$newSheet.Cells.Item(14,4).Value = 8
$newSheet.Cells.Item(14,5).Value = 7
$newSheet.Cells.Item(14,6).Value = 0
$newSheet.Cells.Item(15,1).Value = 'ALGORITHM_STATUS'
$newSheet.Cells.Item(15,2).Value = '0X06'
$newSheet.Cells.Item(15,3).Value = '0X00'
$newSheet.Cells.Item(15,4).Value = 8
$newSheet.Cells.Item(15,5).Value = 7
$newSheet.Cells.Item(15,6).Value = 0
$newSheet.Cells.Item(16,1).Value = 'CALIBRATION_STATUS'
$newSheet.Cells.Item(16,2).Value = '0X07'
$newSheet.Cells.Item(16,3).Value = '0X00'
$newSheet.Cells.Item(16,4).Value = 8
$newSheet.Cells.Item(16,5).Value = 7
$newSheet.Cells.Item(16,6).Value = 0
$newSheet.Cells.Item(17,1).Value = 'CMD_STAT'
$newSheet.Cells.Item(17,2).Value = '0X08'
$newSheet.Cells.Item(17,3).Value = '0X00'
$newSheet.Cells.Item(17,4).Value = 8
$newSheet.Cells.Item(17,5).Value = 7
$newSheet.Cells.Item(17,6).Value = 0
$newSheet.Cells.Item(18,1).Value = 'PREV_CMD'
$newSheet.Cells.Item(18,2).Value = '0X09'
$newSheet.Cells.Item(18,3).Value = '0X00'
$newSheet.Cells.Item(18,4).Value = 8
$newSheet.Cells.Item(18,5).Value = 7
$newSheet.Cells.Item(18,6).Value = 0
$newSheet.Cells.Item(19,1).Value = 'MODE'
$newSheet.Cells.Item(19,2).Value = '0X10'
$newSheet.Cells.Item(19,3).Value = '0X00'
$newSheet.Cells.Item(19,4).Value = 8
$newSheet.Cells.Item(19,5).Value = 7
$newSheet.Cells.Item(19,6).Value = 0
$newSheet.Cells.Item(20,1).Value = 'LIVE_BEAT'
$newSheet.Cells.Item(20,2).Value = '0X0A'
$newSheet.Cells.Item(20,3).Value = '0X00'
$newSheet.Cells.Item(20,4).Value = 8
$newSheet.Cells.Item(20,5).Value = 7
$newSheet.Cells.Item(20,6).Value = 0
$newSheet.Cells.Item(21,1).Value = 'ACTIVE_RANGE'
$newSheet.Cells.Item(21,2).Value = '0X19'
$newSheet.Cells.Item(21,3).Value = '0X6F'
$newSheet.Cells.Item(21,4).Value = 8
$newSheet.Cells.Item(21,5).Value = 7
$newSheet.Cells.Item(21,6).Value = 0
$newSheet.Cells.Item(22,1).Value = 'CONFIG_RESULT'
$newSheet.Cells.Item(22,2).Value = '0X20'
$newSheet.Cells.Item(22,3).Value = '0X00'
$newSheet.Cells.Item(22,4).Value = 8
$newSheet.Cells.Item(22,5).Value = 7
$newSheet.Cells.Item(22,6).Value = 0
$newSheet.Cells.Item(23,1).Value = 'TID'
$newSheet.Cells.Item(23,2).Value = '0X21'
$newSheet.Cells.Item(23,3).Value = '0X00'
$newSheet.Cells.Item(23,4).Value = 8
$newSheet.Cells.Item(23,5).Value = 7
$newSheet.Cells.Item(23,6).Value = 0
$newSheet.Cells.Item(24,1).Value = 'SIZE'
$newSheet.Cells.Item(24,2).Value = '0X22'
$newSheet.Cells.Item(24,3).Value = '0X00'
$newSheet.Cells.Item(24,4).Value = 16
$newSheet.Cells.Item(24,5).Value = 15
$newSheet.Cells.Item(24,6).Value = 0
$newSheet.Cells.Item(25,1).Value = 'PERIOD'
$newSheet.Cells.Item(25,2).Value = '0X24'
$newSheet.Cells.Item(25,3).Value = '0X21'
$newSheet.Cells.Item(25,4).Value = 16
$newSheet.Cells.Item(25,5).Value = 15
$newSheet.Cells.Item(25,6).Value = 0
$newSheet.Cells.Item(26,1).Value = 'KILO_ITERATIONS'
$newSheet.Cells.Item(26,2).Value = '0X26'
$newSheet.Cells.Item(26,3).Value = '0X219'
$newSheet.Cells.Item(26,4).Value = 16
$newSheet.Cells.Item(26,5).Value = 15
$newSheet.Cells.Item(26,6).Value = 0
$newSheet.Cells.Item(27,1).Value = 'CONFIDENCE_THRESHOLD'
$newSheet.Cells.Item(27,2).Value = '0X30'
$newSheet.Cells.Item(27,3).Value = '0X06'
$newSheet.Cells.Item(27,4).Value = 8
$newSheet.Cells.Item(27,5).Value = 7
$newSheet.Cells.Item(27,6).Value = 0
$newSheet.Cells.Item(28,1).Value = 'SPAD_MAP_ID'
$newSheet.Cells.Item(28,2).Value = '0x34'
$newSheet.Cells.Item(28,3).Value = '0x1'
$newSheet.Cells.Item(28,4).Value = 4
$newSheet.Cells.Item(28,5).Value = 3
$newSheet.Cells.Item(28,6).Value = 0
$newSheet.Cells.Item(29,1).Value = 'ALG_SETTING'
$newSheet.Cells.Item(29,2).Value = '0x35'
$newSheet.Cells.Item(29,3).Value = '0x04'
$newSheet.Cells.Item(29,4).Value = 8
$newSheet.Cells.Item(29,5).Value = 7
$newSheet.Cells.Item(29,6).Value = 0
$newSheet.Cells.Item(30,1).Value = 'HIST_DUMP'
$newSheet.Cells.Item(30,2).Value = '0x39'
$newSheet.Cells.Item(30,3).Value = '0x00'
$newSheet.Cells.Item(30,4).Value = 8
$newSheet.Cells.Item(30,5).Value = 7
$newSheet.Cells.Item(30,6).Value = 0
$newSheet.Cells.Item(31,1).Value = 'SPREAD_SPECTRUM'
$newSheet.Cells.Item(31,2).Value = '0X3A'
$newSheet.Cells.Item(31,3).Value = '0X00'
$newSheet.Cells.Item(31,4).Value = 3
$newSheet.Cells.Item(31,5).Value = 2
$newSheet.Cells.Item(31,6).Value = 0
$newSheet.Cells.Item(32,1).Value = 'I2C_SLAVE_ADDRESS'
$newSheet.Cells.Item(32,2).Value = '0X3B'
$newSheet.Cells.Item(32,3).Value = '0X41'
$newSheet.Cells.Item(32,4).Value = 7
$newSheet.Cells.Item(32,5).Value = 7
$newSheet.Cells.Item(32,6).Value = 1
$newSheet.Cells.Item(33,1).Value = 'OSC_TRIM_VALUE'
$newSheet.Cells.Item(33,2).Value = '0X3C'
$newSheet.Cells.Item(33,3).Value = '0X00'
$newSheet.Cells.Item(33,4).Value = 9
$newSheet.Cells.Item(33,5).Value = 8
$newSheet.Cells.Item(33,6).Value = 0
$newSheet.Cells.Item(34,1).Value = 'I2C_ADDR_CHANGE'
$newSheet.Cells.Item(34,2).Value = '0X3E'
$newSheet.Cells.Item(34,3).Value = '0X00'
$newSheet.Cells.Item(34,4).Value = 8
$newSheet.Cells.Item(34,5).Value = 7
$newSheet.Cells.Item(34,6).Value = 0
$newSheet.Cells.Item(35,1).Value = '                                                                                                                                                  '

# ------------------------------------------------------------------
# 3. Column widths / formatting to roughly match the authored sheet.
# ------------------------------------------------------------------
$newSheet.Columns.Item(1).ColumnWidth = 24
$newSheet.Columns.Item(2).ColumnWidth = 12
$newSheet.Columns.Item(3).ColumnWidth = 13
$newSheet.Columns.Item(4).ColumnWidth = 9
$newSheet.Columns.Item(5).ColumnWidth = 15
$newSheet.Columns.Item(6).ColumnWidth = 14

# ------------------------------------------------------------------
# 4. Remove the duplicated "8828..." rows that used to be appended to
#    the bottom of the TMF8801 sheet (rows 43-46): AKG_SETTINGS,
#    HIST_DUMP, ACTIVE_RANGE and SPAD_MAP_ID. Those registers now live
#    on the dedicated TMF8828 sheet created above.
# ------------------------------------------------------------------
$tmf8801.Range("A43:F46").EntireRow.Delete() | Out-Null

# ------------------------------------------------------------------
# 5. Restore a sensible view/selection state on both sheets and make
#    TMF8828 the active tab (mirrors the authored workbook).
# ------------------------------------------------------------------
$tmf8801.Activate()
$tmf8801.Range("E44").Select() | Out-Null

$newSheet.Activate()
$newSheet.Range("A35").Select() | Out-Null
